$d = $word.ActiveDocument

# NOTE: in this headless runtime, Range.Find.Execute with a Replace action
# (wdReplaceOne/wdReplaceAll) matches and replaces against the *whole*
# document body, ignoring the bounds of the Range/Cell Find was invoked on.
# Several answer cells in this worksheet table transiently or originally
# share identical text, so Find-based replacement would edit the wrong cell.
# Assigning Range.Text directly is properly scoped to the target Range, so
# every value below is updated that way instead.

# Update the date line at the top of the document
$d.Paragraphs.Item(1).Range.Text = "2023-09-30 Saturday"

# Update every arithmetic answer in the single results table, cell by cell
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "50-9=41"
$t.Cell(1, 2).Range.Text = "72-29=43"
$t.Cell(1, 3).Range.Text = "46+36=82"
$t.Cell(1, 4).Range.Text = "81-33=48"
$t.Cell(1, 5).Range.Text = "58+23=81"

$t.Cell(2, 1).Range.Text = "82-18=64"
$t.Cell(2, 2).Range.Text = "27+64=91"
$t.Cell(2, 3).Range.Text = "88+5=93"
$t.Cell(2, 4).Range.Text = "70-62=8"
$t.Cell(2, 5).Range.Text = "84-48=36"

$t.Cell(3, 1).Range.Text = "51-37=14"
$t.Cell(3, 2).Range.Text = "18+59=77"
$t.Cell(3, 3).Range.Text = "50-9=41"
$t.Cell(3, 4).Range.Text = "58-19=39"
$t.Cell(3, 5).Range.Text = "75-8=67"

$t.Cell(4, 1).Range.Text = "96-18=78"
$t.Cell(4, 2).Range.Text = "62-49=13"
$t.Cell(4, 3).Range.Text = "40-19=21"
$t.Cell(4, 4).Range.Text = "49+19=68"
$t.Cell(4, 5).Range.Text = "25-7=18"

$t.Cell(5, 1).Range.Text = "63+28=91"
$t.Cell(5, 2).Range.Text = "84-37=47"
$t.Cell(5, 3).Range.Text = "3+68=71"
$t.Cell(5, 4).Range.Text = "5+87=92"
$t.Cell(5, 5).Range.Text = "44-9=35"

$t.Cell(6, 1).Range.Text = "7+25=32"
$t.Cell(6, 2).Range.Text = "6+29=35"
$t.Cell(6, 3).Range.Text = "23+59=82"
$t.Cell(6, 4).Range.Text = "97-19=78"
$t.Cell(6, 5).Range.Text = "98-29=69"

$t.Cell(7, 1).Range.Text = "63-59=4"
$t.Cell(7, 2).Range.Text = "48+39=87"
$t.Cell(7, 3).Range.Text = "55+17=72"
$t.Cell(7, 4).Range.Text = "82-46=36"
$t.Cell(7, 5).Range.Text = "64-49=15"

$t.Cell(8, 1).Range.Text = "44+18=62"
$t.Cell(8, 2).Range.Text = "85-27=58"
$t.Cell(8, 3).Range.Text = "86+8=94"
$t.Cell(8, 4).Range.Text = "20-4=16"
$t.Cell(8, 5).Range.Text = "69+4=73"

$t.Cell(9, 1).Range.Text = "58+7=65"
$t.Cell(9, 2).Range.Text = "76-7=69"
$t.Cell(9, 3).Range.Text = "85-6=79"
$t.Cell(9, 4).Range.Text = "46+9=55"
$t.Cell(9, 5).Range.Text = "93-25=68"

$t.Cell(10, 1).Range.Text = "63-19=44"
$t.Cell(10, 2).Range.Text = "64+27=91"
$t.Cell(10, 3).Range.Text = "76-18=58"
$t.Cell(10, 4).Range.Text = "74+19=93"
$t.Cell(10, 5).Range.Text = "70-55=15"

$t.Cell(11, 1).Range.Text = "84-15=69"
$t.Cell(11, 2).Range.Text = "19+64=83"
$t.Cell(11, 3).Range.Text = "6+29=35"
$t.Cell(11, 4).Range.Text = "74+19=93"
$t.Cell(11, 5).Range.Text = "36-28=8"

$t.Cell(12, 1).Range.Text = "26+18=44"
$t.Cell(12, 2).Range.Text = "30-12=18"
$t.Cell(12, 3).Range.Text = "28+9=37"
$t.Cell(12, 4).Range.Text = "70-54=16"
$t.Cell(12, 5).Range.Text = "47+5=52"

$t.Cell(13, 1).Range.Text = "17+38=55"
$t.Cell(13, 2).Range.Text = "95-76=19"
$t.Cell(13, 3).Range.Text = "46+47=93"
$t.Cell(13, 4).Range.Text = "56+17=73"
$t.Cell(13, 5).Range.Text = "80-59=21"

$t.Cell(14, 1).Range.Text = "33-29=4"
$t.Cell(14, 2).Range.Text = "53-4=49"
$t.Cell(14, 3).Range.Text = "55-38=17"
$t.Cell(14, 4).Range.Text = "70-52=18"
$t.Cell(14, 5).Range.Text = "74-28=46"

$t.Cell(15, 1).Range.Text = "59+38=97"
$t.Cell(15, 2).Range.Text = "44+8=52"
$t.Cell(15, 3).Range.Text = "6+8=14"
$t.Cell(15, 4).Range.Text = "8+37=45"
$t.Cell(15, 5).Range.Text = "92-46=46"

$t.Cell(16, 1).Range.Text = "73-64=9"
$t.Cell(16, 2).Range.Text = "96-27=69"
$t.Cell(16, 3).Range.Text = "27-8=19"
$t.Cell(16, 4).Range.Text = "90-79=11"
$t.Cell(16, 5).Range.Text = "47+48=95"

$t.Cell(17, 1).Range.Text = "16+7=23"
$t.Cell(17, 2).Range.Text = "73-59=14"
$t.Cell(17, 3).Range.Text = "75-16=59"
$t.Cell(17, 4).Range.Text = "39+44=83"
$t.Cell(17, 5).Range.Text = "48+23=71"

$t.Cell(18, 1).Range.Text = "39+22=61"
$t.Cell(18, 2).Range.Text = "88+5=93"
$t.Cell(18, 3).Range.Text = "58+8=66"
$t.Cell(18, 4).Range.Text = "78+18=96"
$t.Cell(18, 5).Range.Text = "92-35=57"

$t.Cell(19, 1).Range.Text = "27+39=66"
$t.Cell(19, 2).Range.Text = "54+27=81"
$t.Cell(19, 3).Range.Text = "94-17=77"
$t.Cell(19, 4).Range.Text = "86-17=69"
$t.Cell(19, 5).Range.Text = "89+6=95"

$t.Cell(20, 1).Range.Text = "38+29=67"
$t.Cell(20, 2).Range.Text = "94-76=18"
$t.Cell(20, 3).Range.Text = "9+17=26"
$t.Cell(20, 4).Range.Text = "32-17=15"
$t.Cell(20, 5).Range.Text = "81-28=53"

